# Apply the commit's changes to the assets/liabilities workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "Summary": update name, income, totals, net worth, ratio ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Khalid Al Kalbani"
$wsSummary.Range("B4").Value = 2538.86
$wsSummary.Range("B6").Value = 4857
$wsSummary.Range("B7").Value = 29098
$wsSummary.Range("B8").Value = -24241
$wsSummary.Range("B9").Value = 0.17

# --- Sheet "Assets": remove the two "Vehicles" line items, keep only
#     "Liquid Assets / Savings Account" and "TOTAL ASSETS", update values ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Rows("2:3").Delete()
$wsAssets.Range("C2").Value = 4857
$wsAssets.Range("C3").Value = 4857

# --- Sheet "Liabilities": remove the Auto Loan x2 and Personal Loan line
#     items, keep only "Credit Cards / Credit Card Balance" and
#     "TOTAL LIABILITIES", update values ---
$wsLiabilities = $wb.Worksheets.Item("Liabilities")
$wsLiabilities.Rows("2:4").Delete()
$wsLiabilities.Range("C2").Value = 29098
$wsLiabilities.Range("D2").Value = 1455
$wsLiabilities.Range("E2").Value = 1
$wsLiabilities.Range("C3").Value = 29098
